$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_4_6_0"
$ws.Range("B2").Value = 0.5013757786130772
$ws.Range("C2").Value = 0.4028849805350341
$ws.Range("D2").Value = 0.9054571695119211
$ws.Range("E2").Value = 0.6190784619707741
$ws.Range("F2").Value = 0.5518300533294678
$ws.Range("G2").Value = 0.738880455493927
$ws.Range("H2").Value = 0.07864782214164734
$ws.Range("I2").Value = 0.4281877279281616

$ws.Range("A3").Value = "model_4_6_1"
$ws.Range("B3").Value = 0.8137628730313161
$ws.Range("C3").Value = 0.8775240629631477
$ws.Range("D3").Value = 0.3758037821782957
$ws.Range("E3").Value = 0.7112415491375191
$ws.Range("F3").Value = 0.2061096280813217
$ws.Range("G3").Value = 0.1515538543462753
$ws.Range("H3").Value = 0.5192532539367676
$ws.Range("I3").Value = 0.3245887160301208

$ws.Range("A4").Value = "model_4_6_3"
$ws.Range("B4").Value = 0.8232035498096431
$ws.Range("C4").Value = 0.8328014063406032
$ws.Range("D4").Value = 0.2922395505951255
$ws.Range("E4").Value = 0.656075584589394
$ws.Range("F4").Value = 0.1956615746021271
$ws.Range("G4").Value = 0.2068944573402405
$ws.Range("H4").Value = 0.5887681841850281
$ws.Range("I4").Value = 0.3865998685359955

$ws.Range("A5").Value = "model_4_6_2"
$ws.Range("B5").Value = 0.824145197997904
$ws.Range("C5").Value = 0.8711611522063725
$ws.Range("D5").Value = 0.2745986205837299
$ws.Range("E5").Value = 0.672286815611612
$ws.Range("F5").Value = 0.194619432091713
$ws.Range("G5").Value = 0.1594274342060089
$ws.Range("H5").Value = 0.6034433245658875
$ws.Range("I5").Value = 0.3683770596981049

$ws.Range("A6").Value = "model_4_6_4"
$ws.Range("B6").Value = 0.8275009426238287
$ws.Range("C6").Value = 0.8378948554128983
$ws.Range("D6").Value = 0.3104568055056217
$ws.Range("E6").Value = 0.6653887441676474
$ws.Range("F6").Value = 0.1909056454896927
$ws.Range("G6").Value = 0.2005917280912399
$ws.Range("H6").Value = 0.5736137628555298
$ws.Range("I6").Value = 0.3761310875415802

$ws.Range("A7").Value = "model_4_6_5"
$ws.Range("B7").Value = 0.8314333949474141
$ws.Range("C7").Value = 0.8344704463030714
$ws.Range("D7").Value = 0.3492421212644471
$ws.Range("E7").Value = 0.6768999783048337
$ws.Range("F7").Value = 0.186553567647934
$ws.Range("G7").Value = 0.204829141497612
$ws.Range("H7").Value = 0.5413492321968079
$ws.Range("I7").Value = 0.3631915152072906

$ws.Range("A8").Value = "model_4_6_6"
$ws.Range("B8").Value = 0.8353511062756684
$ws.Range("C8").Value = 0.8364547613302452
$ws.Range("D8").Value = 0.3760489314268468
$ws.Range("E8").Value = 0.6873920476892962
$ws.Range("F8").Value = 0.1822178065776825
$ws.Range("G8").Value = 0.2023737132549286
$ws.Range("H8").Value = 0.519049346446991
$ws.Range("I8").Value = 0.3513975441455841

$ws.Range("A9").Value = "model_4_6_7"
$ws.Range("B9").Value = 0.8387856167032665
$ws.Range("C9").Value = 0.8402038965031353
$ws.Range("D9").Value = 0.3961585964543227
$ws.Range("E9").Value = 0.6965796863135352
$ws.Range("F9").Value = 0.178416833281517
$ws.Range("G9").Value = 0.197734460234642
$ws.Range("H9").Value = 0.5023205280303955
$ws.Range("I9").Value = 0.3410698771476746

$ws.Range("A10").Value = "model_4_6_8"
$ws.Range("B10").Value = 0.8402185903602539
$ws.Range("C10").Value = 0.848357712776359
$ws.Range("D10").Value = 0.3866037855174403
$ws.Range("E10").Value = 0.6980047751341572
$ws.Range("F10").Value = 0.1768309324979782
$ws.Range("G10").Value = 0.1876447945833206
$ws.Range("H10").Value = 0.5102689862251282
$ws.Range("I10").Value = 0.3394679427146912

$ws.Range("A11").Value = "model_4_6_9"
$ws.Range("B11").Value = 0.8431099315661166
$ws.Range("C11").Value = 0.8523317874685616
$ws.Range("D11").Value = 0.3968175733680811
$ws.Range("E11").Value = 0.7038780199257281
$ws.Range("F11").Value = 0.1736310720443726
$ws.Range("G11").Value = 0.1827272027730942
$ws.Range("H11").Value = 0.501772403717041
$ws.Range("I11").Value = 0.3328659534454346

$ws.Range("A12").Value = "model_4_6_10"
$ws.Range("B12").Value = 0.8449530122960565
$ws.Range("C12").Value = 0.8537026758499656
$ws.Range("D12").Value = 0.4040686850507513
$ws.Range("E12").Value = 0.7072019002672183
$ws.Range("F12").Value = 0.1715913116931915
$ws.Range("G12").Value = 0.1810308396816254
$ws.Range("H12").Value = 0.4957403838634491
$ws.Range("I12").Value = 0.3291296064853668

$ws.Range("A13").Value = "model_4_6_11"
$ws.Range("B13").Value = 0.8464310241482533
$ws.Range("C13").Value = 0.8632087012668473
$ws.Range("D13").Value = 0.3841772660525666
$ws.Range("E13").Value = 0.7058145025202875
$ws.Range("F13").Value = 0.1699555963277817
$ws.Range("G13").Value = 0.1692679226398468
$ws.Range("H13").Value = 0.5122875571250916
$ws.Range("I13").Value = 0.3306891620159149

$ws.Range("A14").Value = "model_4_6_24"
$ws.Range("B14").Value = 0.8467292770655194
$ws.Range("C14").Value = 0.7418947016777261
$ws.Range("D14").Value = 0.5003045194482674
$ws.Range("E14").Value = 0.6755562587778342
$ws.Range("F14").Value = 0.169625535607338
$ws.Range("G14").Value = 0.319383978843689
$ws.Range("H14").Value = 0.4156841933727264
$ws.Range("I14").Value = 0.3647019863128662

$ws.Range("A15").Value = "model_4_6_23"
$ws.Range("B15").Value = 0.8486743129637274
$ws.Range("C15").Value = 0.7520328972847349
$ws.Range("D15").Value = 0.5033619167712369
$ws.Range("E15").Value = 0.682529543179228
$ws.Range("F15").Value = 0.1674729436635971
$ws.Range("G15").Value = 0.3068387806415558
$ws.Range("H15").Value = 0.4131408035755157
$ws.Range("I15").Value = 0.356863409280777

$ws.Range("A16").Value = "model_4_6_22"
$ws.Range("B16").Value = 0.8515166852667585
$ws.Range("C16").Value = 0.7597670136164008
$ws.Range("D16").Value = 0.5231190291556438
$ws.Range("E16").Value = 0.6939169274880409
$ws.Range("F16").Value = 0.1643272787332535
$ws.Range("G16").Value = 0.2972684800624847
$ws.Range("H16").Value = 0.3967053890228271
$ws.Range("I16").Value = 0.3440630435943604

$ws.Range("A17").Value = "model_4_6_21"
$ws.Range("B17").Value = 0.8570999737400755
$ws.Range("C17").Value = 0.782114018340464
$ws.Range("D17").Value = 0.5429528054707176
$ws.Range("E17").Value = 0.7138484476085311
$ws.Range("F17").Value = 0.1581482142210007
$ws.Range("G17").Value = 0.2696158885955811
$ws.Range("H17").Value = 0.3802061378955841
$ws.Range("I17").Value = 0.3216583132743835

$ws.Range("A18").Value = "model_4_6_20"
$ws.Range("B18").Value = 0.8586633426390097
$ws.Range("C18").Value = 0.7994783336245628
$ws.Range("D18").Value = 0.5374304382425685
$ws.Range("E18").Value = 0.72204477817114
$ws.Range("F18").Value = 0.1564180105924606
$ws.Range("G18").Value = 0.2481289952993393
$ws.Range("H18").Value = 0.3848000466823578
$ws.Range("I18").Value = 0.3124449551105499

$ws.Range("A19").Value = "model_4_6_19"
$ws.Range("B19").Value = 0.8619858479808258
$ws.Range("C19").Value = 0.8082670097948172
$ws.Range("D19").Value = 0.5592723231629828
$ws.Range("E19").Value = 0.7347734088569864
$ws.Range("F19").Value = 0.1527410000562668
$ws.Range("G19").Value = 0.237253725528717
$ws.Range("H19").Value = 0.3666303157806396
$ws.Range("I19").Value = 0.2981369197368622

$ws.Range("A20").Value = "model_4_6_18"
$ws.Range("B20").Value = 0.8629205629763658
$ws.Range("C20").Value = 0.81442893653408
$ws.Range("D20").Value = 0.5640853171565566
$ws.Range("E20").Value = 0.7400408153354606
$ws.Range("F20").Value = 0.1517065465450287
$ws.Range("G20").Value = 0.22962886095047
$ws.Range("H20").Value = 0.3626265525817871
$ws.Range("I20").Value = 0.2922159135341644

$ws.Range("A21").Value = "model_4_6_17"
$ws.Range("B21").Value = 0.8630385039166319
$ws.Range("C21").Value = 0.8159382878378159
$ws.Range("D21").Value = 0.5641476916585437
$ws.Range("E21").Value = 0.7409421875002662
$ws.Range("F21").Value = 0.1515760272741318
$ws.Range("G21").Value = 0.2277611643075943
$ws.Range("H21").Value = 0.3625746369361877
$ws.Range("I21").Value = 0.2912026941776276

$ws.Range("A22").Value = "model_4_6_16"
$ws.Range("B22").Value = 0.8633990572767417
$ws.Range("C22").Value = 0.8183188209938738
$ws.Range("D22").Value = 0.5651298477847846
$ws.Range("E22").Value = 0.7426718024364617
$ws.Range("F22").Value = 0.1511770039796829
$ws.Range("G22").Value = 0.2248154431581497
$ws.Range("H22").Value = 0.3617576062679291
$ws.Range("I22").Value = 0.2892584800720215

$ws.Range("A23").Value = "model_4_6_15"
$ws.Range("B23").Value = 0.8670793598768672
$ws.Range("C23").Value = 0.8258133927761443
$ws.Range("D23").Value = 0.5922304238241984
$ws.Range("E23").Value = 0.7564771232673304
$ws.Range("F23").Value = 0.1471039950847626
$ws.Range("G23").Value = 0.2155415415763855
$ws.Range("H23").Value = 0.3392133116722107
$ws.Range("I23").Value = 0.2737401127815247

$ws.Range("A24").Value = "model_4_6_12"
$ws.Range("B24").Value = 0.8678596145333395
$ws.Range("C24").Value = 0.8542205803309297
$ws.Range("D24").Value = 0.5725861173499572
$ws.Range("E24").Value = 0.766191151626605
$ws.Range("F24").Value = 0.1462404727935791
$ws.Range("G24").Value = 0.1803899854421616
$ws.Range("H24").Value = 0.3555549383163452
$ws.Range("I24").Value = 0.2628207206726074

$ws.Range("A25").Value = "model_4_6_14"
$ws.Range("B25").Value = 0.8700180874251131
$ws.Range("C25").Value = 0.8350790783079319
$ws.Range("D25").Value = 0.6093021058164019
$ws.Range("E25").Value = 0.7678231549003828
$ws.Range("F25").Value = 0.1438516825437546
$ws.Range("G25").Value = 0.2040760070085526
$ws.Range("H25").Value = 0.325011819601059
$ws.Range("I25").Value = 0.2609862387180328

$ws.Range("A26").Value = "model_4_6_13"
$ws.Range("B26").Value = 0.8700368875986981
$ws.Range("C26").Value = 0.843598651617546
$ws.Range("D26").Value = 0.6012481233254092
$ws.Range("E26").Value = 0.7699827340572721
$ws.Range("F26").Value = 0.143830880522728
$ws.Range("G26").Value = 0.1935337334871292
$ws.Range("H26").Value = 0.3317117393016815
$ws.Range("I26").Value = 0.2585587203502655

